# US 39 & 42 and Updated Report
# Checks for anniversaries and reject dates that cannot be converted to a
# datetime object.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Backlog sheet: US39 moves from "Coding" to "Done", and two new backlog
# items (US41 "Include Partial Dates", US42 "Reject Illegitimate Dates")
# are appended.
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

# US39 (row 40) is now complete.
$backlog.Range("E40").Value = "Done"

# New backlog row for US41 - still in progress ("Coding").
$backlog.Range("A42").Value = 41
$backlog.Range("B42").Value = "US41"
$backlog.Range("C42").Value = "Include Partial Dates"
$backlog.Range("D42").Value = "Na"
$backlog.Range("E42").Value = "Coding"
$backlog.Range("B42").Interior.Color = 12444887
$backlog.Range("E42").Interior.Color = 12444887

# New backlog row for US42 - complete ("Done").
$backlog.Range("A43").Value = 42
$backlog.Range("B43").Value = "US42"
$backlog.Range("C43").Value = "Reject Illegitimate Dates"
$backlog.Range("D43").Value = "Na"
$backlog.Range("E43").Value = "Done"

# ---------------------------------------------------------------------
# Sprint4 sheet: close out US39 (anniversaries) with actual size/time,
# and replace the US40 placeholder row with the newly completed US42
# (reject illegitimate / non-datetime dates) story.
# ---------------------------------------------------------------------
$sprint4 = $wb.Worksheets.Item("Sprint4")

# US39 - "List upcoming anniversaries" - completed.
$sprint4.Range("G10").Value = 20
$sprint4.Range("H10").Value = 15
$sprint4.Range("I10").Value = "Done"

# Row 11 becomes US42 "Reject Illegitimate dates" - completed.
$sprint4.Range("A11").Value = "US42"
$sprint4.Range("B11").Value = "Reject Illegitimate dates"
$sprint4.Range("G11").Value = 10
$sprint4.Range("H11").Value = 10
$sprint4.Range("I11").Value = "Done"

# ---------------------------------------------------------------------
# Selection / scroll position updates left behind by the author's last
# editing session in each sheet.
# ---------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Item("Sprint3")
$sprint3.Range("A21").Select()

$sprint4.Range("I10").Select()

# Leave the workbook focused back on the Backlog sheet (the tab that was
# active before/after the edits), with the cursor on the newly added row.
$backlog.Range("E43").Select()
